$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.550.44'
$ws.Range('E2').Value = '  +2.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.350.49'
$ws.Range('E3').Value = '  +5.90%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.43'
$ws.Range('E5').Value = '  +5.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.95'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.645'
$ws.Range('E7').Value = '  +3.28%  '
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.637'
$ws.Range('E9').Value = '  +6.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.26'
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.85'
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('E13').Value = '  +3.43%  '
$ws.Range('E14').Value = '  +2.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.37'
$ws.Range('E15').Value = '  +8.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.705.28'
$ws.Range('E16').Value = '  +5.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.420.09'
$ws.Range('E17').Value = '  +8.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.519.70'
$ws.Range('E18').Value = '  +2.92%  '
$ws.Range('E19').Value = '  +3.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.25'
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '75.73'
$ws.Range('E21').Value = '  +4.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.45'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.57'
$ws.Range('E23').Value = '  +11.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '257.17'
$ws.Range('E24').Value = '  +12.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.14'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.10'
$ws.Range('E26').Value = '  +4.52%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.25'
$ws.Range('E28').Value = '  +3.03%  '
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('E30').Value = '  +7.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.86'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('E33').Value = '  +4.02%  '
$ws.Range('E34').Value = '  +8.14%  '
$ws.Range('E35').Value = '  +5.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.99'
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('E37').Value = '  -3.77%  '
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('E39').Value = '  +2.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.71'
$ws.Range('E40').Value = '  +12.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.47'
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.48'
$ws.Range('E42').Value = '  +14.07%  '
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.82'
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.65'
$ws.Range('E46').Value = '  +4.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.30'
$ws.Range('E47').Value = '  +10.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.30'
$ws.Range('E48').Value = '  +7.86%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.41'
$ws.Range('E51').Value = '  +5.34%  '
